$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21, shifting the existing row 21 (and below) down.
$ws.Rows.Item(21).Insert()

# Populate the new row 21 with the new data entry.
$ws.Cells.Item(21, 1).Value = 7
$ws.Cells.Item(21, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(21, 3).Value = "Ñuble"
$ws.Cells.Item(21, 4).Value = 44615
$ws.Cells.Item(21, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(21, 5).Value = 16
$ws.Cells.Item(21, 6).Value = 100112001
$ws.Cells.Item(21, 7).Value = "Berenjena"
$ws.Cells.Item(21, 8).Value = "Sin especificar"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 100
$ws.Cells.Item(21, 11).Value = 11000
$ws.Cells.Item(21, 12).Value = 12000
$ws.Cells.Item(21, 13).Value = 11500
$ws.Cells.Item(21, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(21, 15).Value = "Región Metropolitana"
$ws.Cells.Item(21, 16).Value = 192
$ws.Cells.Item(21, 17).Value = 60
$ws.Cells.Item(21, 18).Value = "Hortaliza"
